# Update stale "INVAL" test-ID abbreviation to "INPV" (missed changes from PA to AA)
# https://github.com/OWASP/owasp-istg/issues/1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Checklist")

$rows = @(41, 42, 53, 54, 65, 66, 77, 78, 89, 90)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 -replace "INVAL", "INPV"
}
